# Atualização de bases das ligas, do dia: 29-03-2024 às 17:05
# Swap the match-data columns (B:AC) between each pair of rows listed below.
# Column A (row sequence id) is left untouched in every pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsA = @(21, 23, 75, 114, 116, 125, 158, 177)
$rowsB = @(22, 24, 76, 115, 117, 128, 159, 178)

for ($i = 0; $i -lt $rowsA.Count; $i++) {
    $r1 = $rowsA[$i]
    $r2 = $rowsB[$i]

    $addr1 = "B$r1`:AC$r1"
    $addr2 = "B$r2`:AC$r2"

    $range1 = $ws.Range($addr1)
    $range2 = $ws.Range($addr2)

    $tmp = $range1.Value2
    $range1.Value2 = $range2.Value2
    $range2.Value2 = $tmp
}

Write-Output "Swapped $($rowsA.Count) row pairs"
